$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.166777666666666
$ws.Range("H2").Value = 6.500332999999999
$ws.Range("I2").Value = 0.3769884032128669
$ws.Range("J2").Value = 0.376988403212867
$ws.Range("M2").Value = 35.42630833333333
$ws.Range("N2").Value = 106.278925
$ws.Range("O2").Value = 0.3011049743557705
$ws.Range("P2").Value = 0.3011049743557705
$ws.Range("Q2").Value = 76.76093370911387
$ws.Range("R2").Value = 690.8484033820249
$ws.Range("S2").Value = 0.1135130834818332
$ws.Range("T2").Value = 0.1135130834818332

# Row 3
$ws.Range("G3").Value = 2.166777666666666
$ws.Range("H3").Value = 6.500332999999999
$ws.Range("I3").Value = 0.3769884032128669
$ws.Range("J3").Value = 0.376988403212867
$ws.Range("O3").Value = 0.1186257117186547
$ws.Range("P3").Value = 0.1186257117186547
$ws.Range("Q3").Value = 30.24134826372256
$ws.Range("R3").Value = 272.172134373503
$ws.Range("S3").Value = 0.04472051764080551
$ws.Range("T3").Value = 0.04472051764080551

# Row 4
$ws.Range("G4").Value = 2.166777666666666
$ws.Range("H4").Value = 6.500332999999999
$ws.Range("I4").Value = 0.3769884032128669
$ws.Range("J4").Value = 0.376988403212867
$ws.Range("M4").Value = 16.22618433333333
$ws.Range("N4").Value = 48.678553
$ws.Range("O4").Value = 0.1379140262544151
$ws.Range("P4").Value = 0.1379140262544152
$ws.Range("Q4").Value = 35.15853382868322
$ws.Range("R4").Value = 316.426804458149
$ws.Range("S4").Value = 0.05199198853830937
$ws.Range("T4").Value = 0.05199198853830938

# Row 5
$ws.Range("G5").Value = 2.166777666666666
$ws.Range("H5").Value = 6.500332999999999
$ws.Range("I5").Value = 0.3769884032128669
$ws.Range("J5").Value = 0.376988403212867
$ws.Range("M5").Value = 4.270717666666667
$ws.Range("N5").Value = 12.812153
$ws.Range("O5").Value = 0.03629885229369049
$ws.Range("P5").Value = 0.03629885229369049
$ws.Range("Q5").Value = 9.253695660772111
$ws.Range("R5").Value = 83.283260946949
$ws.Range("S5").Value = 0.01368424636465809
$ws.Range("T5").Value = 0.01368424636465809

# Row 6
$ws.Range("G6").Value = 2.166777666666666
$ws.Range("H6").Value = 6.500332999999999
$ws.Range("I6").Value = 0.3769884032128669
$ws.Range("J6").Value = 0.376988403212867
$ws.Range("M6").Value = 7.732288666666666
$ws.Range("N6").Value = 23.196866
$ws.Range("O6").Value = 0.06572038381141176
$ws.Range("P6").Value = 0.06572038381141178
$ws.Range("Q6").Value = 16.75415039515311
$ws.Range("R6").Value = 150.787353556378
$ws.Range("S6").Value = 0.02477582255160087
$ws.Range("T6").Value = 0.02477582255160088

# Row 7
$ws.Range("G7").Value = 2.166777666666666
$ws.Range("H7").Value = 6.500332999999999
$ws.Range("I7").Value = 0.3769884032128669
$ws.Range("J7").Value = 0.376988403212867
$ws.Range("M7").Value = 40.042015
$ws.Range("N7").Value = 120.126045
$ws.Range("O7").Value = 0.3403360515660573
$ws.Range("P7").Value = 0.3403360515660573
$ws.Range("Q7").Value = 86.76214383033165
$ws.Range("R7").Value = 780.859294472985
$ws.Range("S7").Value = 0.1283027446356599
$ws.Range("T7").Value = 0.1283027446356599

# Row 8
$ws.Range("I8").Value = 0.3757968909097267
$ws.Range("J8").Value = 0.3757968909097268
$ws.Range("M8").Value = 35.42630833333333
$ws.Range("N8").Value = 106.278925
$ws.Range("O8").Value = 0.3011049743557705
$ws.Range("P8").Value = 0.3011049743557705
$ws.Range("Q8").Value = 76.51832254087778
$ws.Range("R8").Value = 688.6649028679
$ws.Range("S8").Value = 0.1131543132003516
$ws.Range("T8").Value = 0.1131543132003516

# Row 9
$ws.Range("I9").Value = 0.3757968909097267
$ws.Range("J9").Value = 0.3757968909097268
$ws.Range("O9").Value = 0.1186257117186547
$ws.Range("P9").Value = 0.1186257117186547
$ws.Range("S9").Value = 0.04457917364582397
$ws.Range("T9").Value = 0.04457917364582397

# Row 10
$ws.Range("I10").Value = 0.3757968909097267
$ws.Range("J10").Value = 0.3757968909097268
$ws.Range("M10").Value = 16.22618433333333
$ws.Range("N10").Value = 48.678553
$ws.Range("O10").Value = 0.1379140262544151
$ws.Range("P10").Value = 0.1379140262544152
$ws.Range("Q10").Value = 35.04741150964044
$ws.Range("R10").Value = 315.426703586764
$ws.Range("S10").Value = 0.05182766227925163
$ws.Range("T10").Value = 0.05182766227925165

# Row 11
$ws.Range("I11").Value = 0.3757968909097267
$ws.Range("J11").Value = 0.3757968909097268
$ws.Range("M11").Value = 4.270717666666667
$ws.Range("N11").Value = 12.812153
$ws.Range("O11").Value = 0.03629885229369049
$ws.Range("P11").Value = 0.03629885229369049
$ws.Range("Q11").Value = 9.224448362618224
$ws.Range("R11").Value = 83.020035263564
$ws.Range("S11").Value = 0.01364099583556029
$ws.Range("T11").Value = 0.01364099583556029

# Row 12
$ws.Range("I12").Value = 0.3757968909097267
$ws.Range("J12").Value = 0.3757968909097268
$ws.Range("M12").Value = 7.732288666666666
$ws.Range("N12").Value = 23.196866
$ws.Range("O12").Value = 0.06572038381141176
$ws.Range("P12").Value = 0.06572038381141178
$ws.Range("Q12").Value = 16.70119710493422
$ws.Range("R12").Value = 150.310773944408
$ws.Range("S12").Value = 0.02469751590572248
$ws.Range("T12").Value = 0.02469751590572249

# Row 13
$ws.Range("I13").Value = 0.3757968909097267
$ws.Range("J13").Value = 0.3757968909097268
$ws.Range("M13").Value = 40.042015
$ws.Range("N13").Value = 120.126045
$ws.Range("O13").Value = 0.3403360515660573
$ws.Range("P13").Value = 0.3403360515660573
$ws.Range("Q13").Value = 86.48792276427334
$ws.Range("R13").Value = 778.39130487846
$ws.Range("S13").Value = 0.1278972300430168
$ws.Range("T13").Value = 0.1278972300430168

# Row 14
$ws.Range("G14").Value = 1.420890666666667
$ws.Range("H14").Value = 4.262672
$ws.Range("I14").Value = 0.2472147058774063
$ws.Range("J14").Value = 0.2472147058774063
$ws.Range("M14").Value = 35.42630833333333
$ws.Range("N14").Value = 106.278925
$ws.Range("O14").Value = 0.3011049743557705
$ws.Range("P14").Value = 0.3011049743557705
$ws.Range("Q14").Value = 50.33691086528889
$ws.Range("R14").Value = 453.0321977876
$ws.Range("S14").Value = 0.07443757767358578
$ws.Range("T14").Value = 0.07443757767358579

# Row 15
$ws.Range("G15").Value = 1.420890666666667
$ws.Range("H15").Value = 4.262672
$ws.Range("I15").Value = 0.2472147058774063
$ws.Range("J15").Value = 0.2472147058774063
$ws.Range("O15").Value = 0.1186257117186547
$ws.Range("P15").Value = 0.1186257117186547
$ws.Range("Q15").Value = 19.83112995688356
$ws.Range("R15").Value = 178.480169611952
$ws.Range("S15").Value = 0.02932602043202521
$ws.Range("T15").Value = 0.02932602043202521

# Row 16
$ws.Range("G16").Value = 1.420890666666667
$ws.Range("H16").Value = 4.262672
$ws.Range("I16").Value = 0.2472147058774063
$ws.Range("J16").Value = 0.2472147058774063
$ws.Range("M16").Value = 16.22618433333333
$ws.Range("N16").Value = 48.678553
$ws.Range("O16").Value = 0.1379140262544151
$ws.Range("P16").Value = 0.1379140262544152
$ws.Range("Q16").Value = 23.05563387484622
$ws.Range("R16").Value = 207.500704873616
$ws.Range("S16").Value = 0.03409437543685413
$ws.Range("T16").Value = 0.03409437543685413

# Row 17
$ws.Range("G17").Value = 1.420890666666667
$ws.Range("H17").Value = 4.262672
$ws.Range("I17").Value = 0.2472147058774063
$ws.Range("J17").Value = 0.2472147058774063
$ws.Range("M17").Value = 4.270717666666667
$ws.Range("N17").Value = 12.812153
$ws.Range("O17").Value = 0.03629885229369049
$ws.Range("P17").Value = 0.03629885229369049
$ws.Range("Q17").Value = 6.068222872535112
$ws.Range("R17").Value = 54.61400585281601
$ws.Range("S17").Value = 0.008973610093472108
$ws.Range("T17").Value = 0.00897361009347211

# Row 18
$ws.Range("G18").Value = 1.420890666666667
$ws.Range("H18").Value = 4.262672
$ws.Range("I18").Value = 0.2472147058774063
$ws.Range("J18").Value = 0.2472147058774063
$ws.Range("M18").Value = 7.732288666666666
$ws.Range("N18").Value = 23.196866
$ws.Range("O18").Value = 0.06572038381141176
$ws.Range("P18").Value = 0.06572038381141178
$ws.Range("Q18").Value = 10.98673679843911
$ws.Range("R18").Value = 98.88063118595201
$ws.Range("S18").Value = 0.01624704535408841
$ws.Range("T18").Value = 0.01624704535408842

# Row 19
$ws.Range("G19").Value = 1.420890666666667
$ws.Range("H19").Value = 4.262672
$ws.Range("I19").Value = 0.2472147058774063
$ws.Range("J19").Value = 0.2472147058774063
$ws.Range("M19").Value = 40.042015
$ws.Range("N19").Value = 120.126045
$ws.Range("O19").Value = 0.3403360515660573
$ws.Range("P19").Value = 0.3403360515660573
$ws.Range("Q19").Value = 56.89532538802667
$ws.Range("R19").Value = 512.0579284922401
$ws.Range("S19").Value = 0.08413607688738065
$ws.Range("T19").Value = 0.08413607688738066
